$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp update (13:20 -> 13:50)
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 13:50"

# Re-synced provincial rows: city order shifted within several blocks
# and casos/recuperados/muertes counts updated to the latest snapshot.
# Row 14: Zaragoza
$ws.Range("A14").Value = "Zaragoza"
$ws.Range("B14").Value = 2068
$ws.Range("C14").Value = 259
$ws.Range("D14").Value = 1654
$ws.Range("E14").Value = 155

# Row 15: A Coruña
$ws.Range("A15").Value = "A Coruña"
$ws.Range("B15").Value = 1969
$ws.Range("C15").Value = 333
$ws.Range("D15").Value = 1788
$ws.Range("E15").Value = 67

# Row 16: Albacete
$ws.Range("A16").Value = "Albacete"
$ws.Range("B16").Value = 1933
$ws.Range("C16").Value = 397
$ws.Range("D16").Value = 1678
$ws.Range("E16").Value = 156

# Row 17: Malaga
$ws.Range("A17").Value = "Malaga"
$ws.Range("B17").Value = 1644
$ws.Range("C17").Value = 93
$ws.Range("D17").Value = 1458
$ws.Range("E17").Value = 93

# Row 18: Toledo
$ws.Range("A18").Value = "Toledo"
$ws.Range("B18").Value = 1593
$ws.Range("C18").Value = 397
$ws.Range("D18").Value = 1298
$ws.Range("E18").Value = 205

# Row 19: Pontevedra
$ws.Range("B19").Value = 1536
$ws.Range("C19").Value = 333
$ws.Range("D19").Value = 1411

# Row 22: Sevilla
$ws.Range("A22").Value = "Sevilla"
$ws.Range("B22").Value = 1371
$ws.Range("C22").Value = 20
$ws.Range("D22").Value = 1294
$ws.Range("E22").Value = 57

# Row 23: Salamanca
$ws.Range("A23").Value = "Salamanca"
$ws.Range("B23").Value = 1316
$ws.Range("C23").Value = 235
$ws.Range("D23").Value = 946
$ws.Range("E23").Value = 135

# Row 24: Cantabria
$ws.Range("A24").Value = "Cantabria"
$ws.Range("B24").Value = 1268
$ws.Range("C24").Value = 60
$ws.Range("D24").Value = 1148
$ws.Range("E24").Value = 60

# Row 25: Granada
$ws.Range("A25").Value = "Granada"
$ws.Range("B25").Value = 1230
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 1129
$ws.Range("E25").Value = 86

# Row 26: Caceres
$ws.Range("A26").Value = "Caceres"
$ws.Range("B26").Value = 1212
$ws.Range("C26").Value = 45
$ws.Range("D26").Value = 1012
$ws.Range("E26").Value = 155

# Row 27: Gipuzkoa/Guipuzcoa
$ws.Range("A27").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B27").Value = 1206
$ws.Range("C27").Value = 2165
$ws.Range("D27").Value = 639
$ws.Range("E27").Value = 52

# Row 28: Valladolid
$ws.Range("A28").Value = "Valladolid"
$ws.Range("B28").Value = 1109
$ws.Range("C28").Value = 262
$ws.Range("D28").Value = 758
$ws.Range("E28").Value = 89

# Row 34: Jaen
$ws.Range("A34").Value = "Jaen"
$ws.Range("B34").Value = 788
$ws.Range("C34").Value = 17
$ws.Range("D34").Value = 732
$ws.Range("E34").Value = 39

# Row 35: Cordoba
$ws.Range("A35").Value = "Cordoba"
$ws.Range("B35").Value = 782
$ws.Range("C35").Value = 4
$ws.Range("D35").Value = 753
$ws.Range("E35").Value = 25

# Row 36: Guadalajara
$ws.Range("A36").Value = "Guadalajara"
$ws.Range("B36").Value = 753
$ws.Range("C36").Value = 397
$ws.Range("D36").Value = 618
$ws.Range("E36").Value = 100

# Row 37: Ourense
$ws.Range("A37").Value = "Ourense"
$ws.Range("B37").Value = 751
$ws.Range("C37").Value = 333
$ws.Range("D37").Value = 660
$ws.Range("E37").Value = 22

# Row 38: Cadiz
$ws.Range("A38").Value = "Cadiz"
$ws.Range("B38").Value = 697
$ws.Range("C38").Value = 17
$ws.Range("D38").Value = 661
$ws.Range("E38").Value = 19

# Row 39: Castello/Castellon
$ws.Range("A39").Value = "Castello/Castellon"
$ws.Range("B39").Value = 660
$ws.Range("C39").Value = 9
$ws.Range("D39").Value = 609
$ws.Range("E39").Value = 42

# Row 40: Soria
$ws.Range("A40").Value = "Soria"
$ws.Range("B40").Value = 659
$ws.Range("C40").Value = 90
$ws.Range("D40").Value = 525
$ws.Range("E40").Value = 44

# Row 41: Badajoz
$ws.Range("B41").Value = 625
$ws.Range("C41").Value = 94
$ws.Range("D41").Value = 505
$ws.Range("E41").Value = 26

# Row 42: Lugo
$ws.Range("A42").Value = "Lugo"
$ws.Range("B42").Value = 586
$ws.Range("C42").Value = 333
$ws.Range("D42").Value = 520
$ws.Range("E42").Value = 11

# Row 43: Avila
$ws.Range("A43").Value = "Avila"
$ws.Range("B43").Value = 512
$ws.Range("C43").Value = 132
$ws.Range("D43").Value = 321
$ws.Range("E43").Value = 59

# Row 46: Huesca
$ws.Range("B46").Value = 349
$ws.Range("C46").Value = 35
$ws.Range("D46").Value = 296
$ws.Range("E46").Value = 18

# Row 48: Almeria
$ws.Range("A48").Value = "Almeria"
$ws.Range("B48").Value = 290
$ws.Range("C48").Value = 14
$ws.Range("D48").Value = 258
$ws.Range("E48").Value = 18

# Row 49: Teruel
$ws.Range("A49").Value = "Teruel"
$ws.Range("B49").Value = 283
$ws.Range("C49").Value = 26
$ws.Range("D49").Value = 232
$ws.Range("E49").Value = 25

# Row 51: Huelva
$ws.Range("A51").Value = "Huelva"
$ws.Range("B51").Value = 224
$ws.Range("C51").Value = 2
$ws.Range("D51").Value = 216
$ws.Range("E51").Value = 6

# Row 52: Mallorca
$ws.Range("A52").Value = "Mallorca"
$ws.Range("B52").Value = 210
$ws.Range("C52").Value = 18
$ws.Range("D52").Value = 194
$ws.Range("E52").Value = 12

